$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsLpWOU = $wb.Worksheets.Item("LpWOU")

# Update the unit-conversion note: "cubic km" -> "billion cubic m"
$wsAbout.Range("A9").Value = "For the U.S., the water output unit is billion cubic m, which is equivalent to Tl (teraliters, or 10^12 liters)"

# Update last-selected cell on the "About" sheet, then switch the active
# sheet/selection to "LpWOU" (matches the saved view state in the diff).
$wsAbout.Range("J10").Select()
$wsLpWOU.Activate()
$wsLpWOU.Range("B2").Select()
